$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Additional transactions added to overflow page one" -- append rows 12-19,
# which repeat the same transaction pattern as rows 8-11 (customer 4 / product
# cycle 3,4,1,2 / transaction date 45254 / cost cycle 22,2,12,41), twice over.
$newRows = @(
    @(4, 1, 3, 45254, 22),
    @(4, 1, 4, 45254, 2),
    @(4, 1, 1, 45254, 12),
    @(4, 1, 2, 45254, 41),
    @(4, 1, 3, 45254, 22),
    @(4, 1, 4, 45254, 2),
    @(4, 1, 1, 45254, 12),
    @(4, 1, 2, 45254, 41)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
}

$endRow = $startRow + $newRows.Count - 1

# Reuse the existing "Date" number format (style index already present in the
# workbook, numFmtId 14) for the new column D cells instead of letting a new
# custom number format get created.
$ws.Range("D2").Copy()
$ws.Range("D$startRow`:D$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# PasteSpecial(Formats) can disturb the pasted-in values on some ranges, so
# reassign the date values explicitly afterwards to be safe.
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 4).Value = 45254
}

# Move the active selection to the new last cell entered, as Excel would after
# typing in the final new row.
$ws.Range("A$endRow").Select() | Out-Null
